# Update "Week 15 simulations" target depth data on both the OFF and DEF
# sheets. Row 2 on each sheet corresponds to the "H" (home) split.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 219
$wsOff.Range("C2").Value = 149
$wsOff.Range("D2").Value = 45
$wsOff.Range("E2").Value = 24

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 196
$wsDef.Range("C2").Value = 130
$wsDef.Range("D2").Value = 48
$wsDef.Range("E2").Value = 18
$wsDef.Range("F2").Value = 3
